# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handback DateTime"
# timestamps produced by the handback report generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date for 91b65fd5-... row
$wsOverview.Range("G4").Value = "2016-08-13 04:55:34"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for 91b65fd5-... row
$wsZhCn.Range("H4").Value = "2016-08-13 04:55:26"
$wsZhCn.Range("K4").Value = "2016-08-13 04:55:56"

# de-de sheet: Correspond Handoff Datetime (mirrors Overview G4) and
# Correspond Handback DateTime for 91b65fd5-... row
$wsDeDe.Range("H4").Value = "2016-08-13 04:55:34"
$wsDeDe.Range("K4").Value = "2016-08-13 04:56:09"
